$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.005.47'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '3.504.19'
$ws.Range('E3').Value = '  -2.04%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '201.87'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '551.16'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.02%  '
$ws.Range('D7').Value = '3.497.81'
$ws.Range('E7').Value = '  -2.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.600'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.56%  '
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('E10').Value = '  -3.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '60.77'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +11.83%  '
$ws.Range('E12').Value = '  -4.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000272'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.81'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.80%  '
$ws.Range('D15').Value = '4.076.40'
$ws.Range('E15').Value = '  -1.90%  '
$ws.Range('D16').Value = '3.508.59'
$ws.Range('E16').Value = '  -2.14%  '
$ws.Range('E17').Value = '  -0.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.51'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.25%  '
$ws.Range('D19').Value = '66.753.54'
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.80'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.77%  '
$ws.Range('E21').Value = '  -3.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '387.58'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.01'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.88'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -9.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.43'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.57%  '
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('E27').Value = '  -4.53%  '
$ws.Range('E28').Value = '  -4.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.71'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.84'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.53%  '
$ws.Range('E31').Value = '  -1.89%  '
$ws.Range('E32').Value = '  -9.80%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '682.38'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.74'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '63.01'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.49%  '
$ws.Range('E36').Value = '  -4.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '39.78'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.405'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.72%  '
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range('E40').Value = '  -2.24%  '
$ws.Range('D41').Value = '3.135.26'
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.997'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.27%  '
$ws.Range('E43').Value = '  -3.57%  '
$ws.Range('D44').Value = '0.0₃0705'
$ws.Range('E44').Value = '  -9.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.81'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +17.00%  '
$ws.Range('E46').Value = '  -12.10%  '
$ws.Range('E47').Value = '  +7.80%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0398'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.37%  '
$ws.Range('E49').Value = '  -2.90%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '137.06'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.38%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.31'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.74%  '
